# Adds the "2020" data column (column Q) to the SDG indicators sheet,
# copying the formatting of the existing "2019" column (P) and filling in
# the reported values, then updates the active selection — mirroring the
# manual edit captured in the source XML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> value reported for 2020 (column Q), same order as column P (2019)
$values = [ordered]@{
    4  = 2020    # header row: year label
    5  = 0.02
    6  = 0
    7  = 0
    8  = 0
    9  = 0.54
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $values.Keys) {
    # Clone the 2019 (P) cell's formatting onto the new 2020 (Q) cell first …
    $ws.Range("P$row").Copy()
    $ws.Range("Q$row").PasteSpecial(-4122)  # xlPasteFormats

    # … then write the actual 2020 value.
    $ws.Range("Q$row").Value = $values[$row]
}

$excel.CutCopyMode = $false

# Restore the worksheet's active-cell selection as recorded in the saved file.
$ws.Range("N19").Select()
